$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename worksheets (Login -> Register)
# ---------------------------------------------------------------------------
$wsValid   = $wb.Worksheets.Item(1)
$wsInvalid = $wb.Worksheets.Item(2)

$wsValid.Name   = 'Valid_Register_Test'
$wsInvalid.Name = 'Invalid_Register_Test'

# ---------------------------------------------------------------------------
# 2. Update Invalid_Register_Test (sheet 2) header row
# ---------------------------------------------------------------------------
$wsInvalid.Range('B1').Value = '${firstname}'
$wsInvalid.Range('C1').Value = '${lastname}'
$wsInvalid.Range('D1').Value = '${mailid}'
$wsInvalid.Range('E1').Value = '${password}'
$wsInvalid.Range('F1').Value = '${expected_error}'

# ---------------------------------------------------------------------------
# 3. Update Invalid_Register_Test data rows (2-5)
# ---------------------------------------------------------------------------
$wsInvalid.Range('B2').Value = 'siva'
$wsInvalid.Range('C2').Value = 'balan'
$wsInvalid.Range('D2').Value = 'sivabalan@gmail.com'
$wsInvalid.Range('E2').Value = 'ai_8'
$wsInvalid.Range('F2').Value = ' is too short (minimum is 5 characters)'

$wsInvalid.Range('B3').Value = 'rashmika'
$wsInvalid.Range('C3').Value = 'mohammed'
$wsInvalid.Range('D3').Value = 'rasmoh@gmail.com'
$wsInvalid.Range('E3').Value = 'rt@w'
$wsInvalid.Range('F3').Value = ' is too short (minimum is 5 characters)'

$wsInvalid.Range('B4').Value = 'nirmala'
$wsInvalid.Range('C4').Value = 'raikumar'
$wsInvalid.Range('D4').Value = 'nirmal@hotmail.com'
$wsInvalid.Range('E4').Value = 'wer1'
$wsInvalid.Range('F4').Value = ' is too short (minimum is 5 characters)'

$wsInvalid.Range('B5').Value = 'lekshmi'
$wsInvalid.Range('C5').Value = 'thangam'
$wsInvalid.Range('D5').Value = 'lekshmi@hotmail.com'
$wsInvalid.Range('E5').Value = 'tan'
$wsInvalid.Range('F5').Value = ' is too short (minimum is 5 characters)'

# ---------------------------------------------------------------------------
# 4. Fix up cell styling
#    - B2:B5 no longer carry the "Hyperlink" style (plain first-name text now)
#    - D2:D5 (the mailid column) now carry the "Hyperlink" style
#    - F2 keeps the vertical-centered, size-10 font style
# ---------------------------------------------------------------------------
$wsInvalid.Range('B2:B5').Style = 'Normal'
$wsInvalid.Range('D2:D5').Style = 'Hyperlink'

$wsInvalid.Range('F2').Font.Size = 10
$wsInvalid.Range('F2').VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 5. Rebuild hyperlinks: drop the old mailid (B) links, add new ones on D,
#    plus a new accidental one on E3 ("rt@w" looks like an email/url).
# ---------------------------------------------------------------------------
$wsInvalid.Hyperlinks.Delete()
$wsInvalid.Hyperlinks.Add($wsInvalid.Range('D2'), 'mailto:sivabalan@gmail.com')
$wsInvalid.Hyperlinks.Add($wsInvalid.Range('D3'), 'mailto:rasmoh@gmail.com')
$wsInvalid.Hyperlinks.Add($wsInvalid.Range('D4'), 'mailto:nirmal@hotmail.com')
$wsInvalid.Hyperlinks.Add($wsInvalid.Range('D5'), 'mailto:lekshmi@hotmail.com')
$wsInvalid.Hyperlinks.Add($wsInvalid.Range('E3'), 'mailto:rt@w')

# Re-apply the Hyperlink style cleanly after linking (Hyperlinks.Add stamps
# its own font formatting; normalise back onto the shared "Hyperlink" style).
$wsInvalid.Range('D2:D5').Style = 'Hyperlink'
$wsInvalid.Range('E3').Style = 'Hyperlink'

# ---------------------------------------------------------------------------
# 6. Column F got wider on the Invalid_Register_Test sheet
# ---------------------------------------------------------------------------
$wsInvalid.Columns.Item(6).ColumnWidth = 35.34

# ---------------------------------------------------------------------------
# 7. Selections: Invalid sheet moves from C8 to F8; Valid sheet moves from
#    E3 to C1:C5. Finish with the Valid sheet active/selected (tab 1) to
#    match the saved workbook state.
# ---------------------------------------------------------------------------
$wsInvalid.Activate()
$wsInvalid.Range('F8').Select()

$wsValid.Activate()
$wsValid.Range('C1:C5').Select()
